# Apply stimulus update:
#  - column L: expand abbreviated response codes to full words
#       r -> right, b -> center, y -> left
#  - columns B, C, D: rename "face" stimulus folder/files to "book"
#       face//face_NN.jpg -> book//book_NN.jpg

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# --- Update column L (correct_ans) abbreviations on every data row ---
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $val = $cell.Value()

    if ($val -eq "r") {
        $cell.Value = "right"
    } elseif ($val -eq "b") {
        $cell.Value = "center"
    } elseif ($val -eq "y") {
        $cell.Value = "left"
    }
}

# --- Rename "face" stimuli to "book" stimuli in columns B, C, D ---
$stimCols = 2, 3, 4   # B, C, D
foreach ($col in $stimCols) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value()

        if ($val -ne $null -and $val.Contains("face//face_")) {
            $cell.Value = $val.Replace("face//face_", "book//book_")
        }
    }
}
